# Applies the "#5: insurance, claim, debt, investment done" edit.
#
# Target sheets: 保險 (Insurance, worksheet #5) and 債務 (Debt, worksheet #6).
# In both sheets the header row had (incorrectly) been filled with a copy
# of row 2's data values instead of generic column names. This edit:
#   1. Rewrites the header row with the correct generic column names
#      (matching the convention used on sheets 土地/建物/存款/基金受益憑證).
#   2. Appends trailing metadata columns (property_category, category,
#      date, legislator_name, legislator_id, source_file, index) to every
#      data row, mirroring the pattern already used on the other sheets.
#
# Cell values are written in natural row-major / left-to-right order so
# that newly introduced shared strings are appended to the shared string
# table in the same order as the canonical edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "保險" (Insurance) -> worksheet #5
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Pre-format the new "date" column as Text so the "2012-04-30" strings
# written below are stored as literal text rather than being
# auto-converted to a date serial number.
$ws5.Range("G2:G5").NumberFormat = "@"

# --- Row 1: header ---
$ws5.Range("B1").Value = "company"
$ws5.Range("C1").Value = "name"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "property_category"
$ws5.Range("F1").Value = "category"
$ws5.Range("G1").Value = "date"
$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("J1").Value = "source_file"
$ws5.Range("K1").Value = "index"

# Copy the bold / bordered header style from B1 onto the new header cells.
$ws5.Range("B1").Copy()
$ws5.Range("E1:K1").PasteSpecial(-4122)

# --- Rows 2-5: append new metadata columns (B/C/D already hold the
# correct data and are left untouched) ---
$ws5.Range("E2").Value = "insurance"
$ws5.Range("F2").Value = "normal"
$ws5.Range("G2").Value = "2012-04-30"
$ws5.Range("H2").Value = "洪秀柱"
$ws5.Range("I2").Value = 546
$ws5.Range("J2").Value = "tmp31791"
$ws5.Range("K2").Value = 81

$ws5.Range("E3").Value = "insurance"
$ws5.Range("F3").Value = "normal"
$ws5.Range("G3").Value = "2012-04-30"
$ws5.Range("H3").Value = "洪秀柱"
$ws5.Range("I3").Value = 546
$ws5.Range("J3").Value = "tmp31791"
$ws5.Range("K3").Value = 82

$ws5.Range("E4").Value = "insurance"
$ws5.Range("F4").Value = "normal"
$ws5.Range("G4").Value = "2012-04-30"
$ws5.Range("H4").Value = "洪秀柱"
$ws5.Range("I4").Value = 546
$ws5.Range("J4").Value = "tmp31791"
$ws5.Range("K4").Value = 83

$ws5.Range("E5").Value = "insurance"
$ws5.Range("F5").Value = "normal"
$ws5.Range("G5").Value = "2012-04-30"
$ws5.Range("H5").Value = "洪秀柱"
$ws5.Range("I5").Value = 546
$ws5.Range("J5").Value = "tmp31791"
$ws5.Range("K5").Value = 84

# Re-apply the plain data-row style (matching column B/C/D) on top of the
# new E:K cells so the Text-forced number format used above doesn't
# linger on the saved cells.
$ws5.Range("B2").Copy()
$ws5.Range("E2:K5").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "債務" (Debt) -> worksheet #6
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Pre-format the new "date" column as Text (same reasoning as above).
$ws6.Range("J2").NumberFormat = "@"

# --- Row 1: header ---
$ws6.Range("B1").Value = "species"
$ws6.Range("C1").Value = "debtor"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

# Copy the bold / bordered header style from B1 onto the new header cells.
$ws6.Range("B1").Copy()
$ws6.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2: append new metadata columns (B/C/D/E/F/G already hold the
# correct data and are left untouched) ---
$ws6.Range("H2").Value = "debt"
$ws6.Range("I2").Value = "normal"
$ws6.Range("J2").Value = "2012-04-30"
$ws6.Range("K2").Value = "洪秀柱"
$ws6.Range("L2").Value = 546
$ws6.Range("M2").Value = "tmp31791"
$ws6.Range("N2").Value = 94

# Re-apply the plain data-row style on top of H2:N2.
$ws6.Range("C2").Copy()
$ws6.Range("H2:N2").PasteSpecial(-4122)
